$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44435
$ws.Cells.Item(2, 11).Value = 17000
$ws.Cells.Item(2, 12).Value = 18000
$ws.Cells.Item(2, 13).Value = 17500
$ws.Cells.Item(2, 16).Value = 972
$ws.Cells.Item(3, 4).Value = 44435
$ws.Cells.Item(3, 9).Value = "Tercera"
$ws.Cells.Item(3, 10).Value = 120
$ws.Cells.Item(3, 11).Value = 14000
$ws.Cells.Item(3, 12).Value = 15000
$ws.Cells.Item(3, 13).Value = 14500
$ws.Cells.Item(3, 16).Value = 806
$ws.Cells.Item(4, 4).Value = 44412
$ws.Cells.Item(4, 10).Value = 150
$ws.Cells.Item(4, 13).Value = 17500
$ws.Cells.Item(4, 16).Value = 972
$ws.Cells.Item(5, 4).Value = 44776
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(6, 4).Value = 44363
$ws.Cells.Item(6, 10).Value = 140
$ws.Cells.Item(6, 11).Value = 14000
$ws.Cells.Item(6, 12).Value = 15000
$ws.Cells.Item(6, 13).Value = 14500
$ws.Cells.Item(6, 16).Value = 806
$ws.Cells.Item(9, 4).Value = 44391
$ws.Cells.Item(9, 11).Value = 15000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 15500
$ws.Cells.Item(9, 16).Value = 861
$ws.Cells.Item(10, 4).Value = 44533
$ws.Cells.Item(10, 8).Value = "Cultivar XV región"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 6000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 6500
$ws.Cells.Item(10, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 650
$ws.Cells.Item(10, 17).Value = 10
$ws.Cells.Item(11, 4).Value = 44533
$ws.Cells.Item(11, 8).Value = "Cultivar XV región"
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 4000
$ws.Cells.Item(11, 12).Value = 5000
$ws.Cells.Item(11, 13).Value = 4500
$ws.Cells.Item(11, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 450
$ws.Cells.Item(11, 17).Value = 10
$ws.Cells.Item(12, 4).Value = 44757
$ws.Cells.Item(12, 8).Value = "Cultivar XV región"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 150
$ws.Cells.Item(12, 11).Value = 6000
$ws.Cells.Item(12, 12).Value = 6500
$ws.Cells.Item(12, 13).Value = 6250
$ws.Cells.Item(12, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 625
$ws.Cells.Item(12, 17).Value = 10
$ws.Cells.Item(13, 4).Value = 44748
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 250
$ws.Cells.Item(15, 4).Value = 44433
$ws.Cells.Item(15, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 11).Value = 17000
$ws.Cells.Item(15, 12).Value = 18000
$ws.Cells.Item(15, 13).Value = 17500
$ws.Cells.Item(15, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 972
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(16, 4).Value = 44433
$ws.Cells.Item(16, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 9).Value = "Tercera"
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 14000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 14500
$ws.Cells.Item(16, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 16).Value = 806
$ws.Cells.Item(16, 17).Value = 18
$ws.Cells.Item(17, 4).Value = 44211
$ws.Cells.Item(17, 9).Value = "Segunda"
$ws.Cells.Item(17, 10).Value = 140
$ws.Cells.Item(17, 11).Value = 4500
$ws.Cells.Item(17, 12).Value = 5000
$ws.Cells.Item(17, 13).Value = 4750
$ws.Cells.Item(17, 16).Value = 475
$ws.Cells.Item(18, 4).Value = 44762
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 16).Value = 861
$ws.Cells.Item(19, 4).Value = 44526
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 12).Value = 5500
$ws.Cells.Item(19, 13).Value = 5250
$ws.Cells.Item(19, 16).Value = 525
$ws.Cells.Item(20, 4).Value = 44526
$ws.Cells.Item(20, 8).Value = "Cultivar XV región"
$ws.Cells.Item(20, 9).Value = "Segunda"
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 4000
$ws.Cells.Item(20, 12).Value = 4500
$ws.Cells.Item(20, 13).Value = 4250
$ws.Cells.Item(20, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 425
$ws.Cells.Item(20, 17).Value = 10
$ws.Cells.Item(21, 4).Value = 44526
$ws.Cells.Item(21, 8).Value = "Cultivar XV región"
$ws.Cells.Item(21, 9).Value = "Tercera"
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 3000
$ws.Cells.Item(21, 12).Value = 3500
$ws.Cells.Item(21, 13).Value = 3250
$ws.Cells.Item(21, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 325
$ws.Cells.Item(21, 17).Value = 10
$ws.Cells.Item(22, 4).Value = 44377
$ws.Cells.Item(22, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 17000
$ws.Cells.Item(22, 12).Value = 18000
$ws.Cells.Item(22, 13).Value = 17600
$ws.Cells.Item(22, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 16).Value = 978
$ws.Cells.Item(22, 17).Value = 18
$ws.Cells.Item(23, 4).Value = 44554
$ws.Cells.Item(23, 8).Value = "Cultivar XV región"
$ws.Cells.Item(23, 10).Value = 200
$ws.Cells.Item(23, 11).Value = 5000
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 13).Value = 5500
$ws.Cells.Item(23, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 550
$ws.Cells.Item(23, 17).Value = 10
$ws.Cells.Item(24, 4).Value = 44742
$ws.Cells.Item(24, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(24, 10).Value = 250
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 16000
$ws.Cells.Item(24, 13).Value = 15500
$ws.Cells.Item(24, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 16).Value = 861
$ws.Cells.Item(24, 17).Value = 18
$ws.Cells.Item(25, 4).Value = 44454
$ws.Cells.Item(25, 10).Value = 160
$ws.Cells.Item(25, 11).Value = 19000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 13).Value = 19500
$ws.Cells.Item(25, 16).Value = 1083
$ws.Cells.Item(26, 4).Value = 44221
$ws.Cells.Item(26, 10).Value = 140
$ws.Cells.Item(26, 11).Value = 5000
$ws.Cells.Item(26, 12).Value = 6000
$ws.Cells.Item(26, 13).Value = 5500
$ws.Cells.Item(26, 16).Value = 550
$ws.Cells.Item(27, 4).Value = 44769
$ws.Cells.Item(27, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 140
$ws.Cells.Item(27, 11).Value = 17000
$ws.Cells.Item(27, 12).Value = 18000
$ws.Cells.Item(27, 13).Value = 17500
$ws.Cells.Item(27, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(27, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 16).Value = 972
$ws.Cells.Item(27, 17).Value = 18
$ws.Cells.Item(28, 4).Value = 44755
$ws.Cells.Item(28, 11).Value = 17000
$ws.Cells.Item(28, 12).Value = 18000
$ws.Cells.Item(28, 13).Value = 17500
$ws.Cells.Item(28, 16).Value = 972
$ws.Cells.Item(29, 4).Value = 44405
$ws.Cells.Item(29, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 140
$ws.Cells.Item(29, 11).Value = 17000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 17500
$ws.Cells.Item(29, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(29, 16).Value = 972
$ws.Cells.Item(29, 17).Value = 18
$ws.Cells.Item(30, 4).Value = 44783
